# IssuesLog.xlsx edit: "Issue-99 Remove DoubleCurveSegment ... references to pegCounts"
#
# Summary of changes applied:
#  1. Fix typo in issue title (row 35): "Mac manus" -> "Mac menus"
#  2. Add a review comment to H8 ("Is this needed")
#  3. Append two new issue rows (99 and 100) for "Remove DoubleCurbe" / "Remove files that
#     should not be tracked"
#  4. Re-point the selection to the new last cell (G100)
#  5. Reshuffle the historical `_FilterDatabase` defined-name chain the same way the
#     workbook's own autofilter bookkeeping does on every save

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Typo fix -----------------------------------------------------------
$ws.Range("B35").Value2 = "Can JavaFX use the Mac menus"

# --- 2. New comment on row 8 -------------------------------------------------
$ws.Range("H8").Value2 = "Is this needed"

# --- 3. New rows 99 & 100 ----------------------------------------------------
$ws.Range("A99").Value2 = 99
$ws.Range("B99").Value2 = "Remove DoubleCurbe"
$ws.Range("C99").Value2 = -10
$ws.Range("E99").Value2 = 43670
$ws.Range("F99").Value2 = "CLOSED"
$ws.Range("H99").Value2 = "Double curves are not used now that the curves can aggressively eat neighbouring pixels"

$ws.Range("A100").Value2 = 100
$ws.Range("B100").Value2 = "Remove files  that should not be tracked "
$ws.Range("C100").Value2 = -10
$ws.Range("E100").Value2 = 43670
$ws.Range("F100").Value2 = "OPEN"

# --- 4. Update the selected cell --------------------------------------------
$ws.Range("G100").Select()

# --- 5. Defined-name "_FilterDatabase" history shuffle ----------------------
# Rename the deepest suffix first so we never collide with an existing name.
$names = $wb.Names
$names.Item(7).Name = "_xlnm._FilterDatabase_0_0_0_0_0_0"
$names.Item(6).Name = "_xlnm._FilterDatabase_0_0_0_0_0"
$names.Item(5).Name = "_xlnm._FilterDatabase_0_0_0_0"
$names.Item(4).Name = "_xlnm._FilterDatabase_0_0_0"
$names.Item(3).Name = "_xlnm._FilterDatabase_0_0"

# The hidden (real) filter database range and the visible bookkeeping range swap values.
$names.Item(1).RefersTo = "=Sheet1!`$A`$1:`$I`$94"
$names.Item(2).RefersTo = "=Sheet1!`$A`$1:`$I`$98"
